$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '5/5'])"
$ws.Range("A3").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"
$ws.Range("A4").Value = "('Saproling', ['Token Creature — Saproling', '1/1'])"
$ws.Range("A5").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"
$ws.Range("A6").Value = "('Wasp', ['Token Artifact Creature — Insect', 'Flying', '1/1'])"
$ws.Range("A7").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"

$ws.Range("A8:A21").ClearContents() | Out-Null
